# Auto-generated Excel COM-interop edit script
# Applies the cell-level changes described in the commit diff
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---- 展览 ----
$ws1.Range("F2").Value = 870
$ws1.Range("F5").Value = 1179
$ws1.Range("F6").Value = 65
$ws1.Range("F7").Value = 4309
$ws1.Range("F8").Value = 2578
$ws1.Range("F10").Value = 2477
$ws1.Range("F17").Value = 107
$ws1.Range("C19").Value = "杭州·赛马娘only—晴空雏菊"
$ws1.Range("D19").Value = "北干街道萧杭路689号 时尚外滩艺术中心"
$ws1.Range("E19").Value = "2024.04.13 09:00-04.13 18:00"
$ws1.Range("F19").Value = 269
$ws1.Range("G19").Value = 66
$ws1.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=81767"
$ws1.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202402/ViMb8nbw1707122090281.jpeg"
$ws1.Range("B20").NumberFormat = "@"
$ws1.Range("B20").Value = "2024-04-20"
$ws1.Range("C20").Value = "杭州·COMIC WORLD次元创作同人季特典·SP·浙里来消"
$ws1.Range("D20").Value = "德胜东路2539号 梦马汽车小镇"
$ws1.Range("E20").Value = "2024.04.20 10:00-04.21 17:00"
$ws1.Range("F20").Value = 71
$ws1.Range("G20").Value = 58
$ws1.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=82573"
$ws1.Range("I20").Value = "//i2.hdslb.com/bfs/openplatform/202403/UgyVGYJa1709879114323.png"
$ws1.Range("C21").Value = "杭州·EVA ONLY漫展"
$ws1.Range("D21").Value = "文三路199号创业大厦众创空间2层 杭州趣链科技有限公司"
$ws1.Range("E21").Value = "2024.04.20 10:00-04.20 17:00"
$ws1.Range("F21").Value = 18
$ws1.Range("G21").Value = 88
$ws1.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=82988"
$ws1.Range("I21").Value = "//i1.hdslb.com/bfs/openplatform/202403/F9yBq4Qo1710756247458.jpeg"
$ws1.Range("C22").Value = "杭州·SK怀旧展&偶像专场"
$ws1.Range("D22").Value = "沈半路171号 T-Car杭州汽车文化主题公园"
$ws1.Range("E22").Value = "2024.04.20 09:00-04.20 22:00"
$ws1.Range("F22").Value = 468
$ws1.Range("G22").Value = 60
$ws1.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=81764"
$ws1.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202402/mtdbSuTZ1707119415384.jpeg"
$ws1.Range("C23").Value = "杭州·m字刘海少年和粉毛少女only"
$ws1.Range("D23").Value = "康候圣街99号 顺丰创新中心"
$ws1.Range("E23").Value = "2024.04.20 09:00-04.20 17:00"
$ws1.Range("F23").Value = 27
$ws1.Range("G23").Value = 68
$ws1.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=82831"
$ws1.Range("I23").Value = "//i2.hdslb.com/bfs/openplatform/202403/bVvk6Eky1710383662942.jpeg"
$ws1.Range("C24").Value = "杭州·【海潮的回响Echo of The Tide】 | 刀客塔们的大群融入派对·明日方舟SPECIAL ONLY"
$ws1.Range("D24").Value = "保淑路2号 The Queen皇后"
$ws1.Range("E24").Value = "2024.04.20 14:00-04.20 18:00"
$ws1.Range("F24").Value = 86
$ws1.Range("G24").Value = 139
$ws1.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=82068"
$ws1.Range("I24").Value = "//i0.hdslb.com/bfs/openplatform/202402/gAR8Svsc1708919248994.png"
$ws1.Range("C25").Value = "杭州·白日梦次元动漫嘉年华"
$ws1.Range("D25").Value = "黄姑山路51-4号 0101park"
$ws1.Range("E25").Value = "2024.04.20 10:00-04.21 18:00"
$ws1.Range("F25").Value = 523
$ws1.Range("G25").Value = 68
$ws1.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=81634"
$ws1.Range("I25").Value = "//i1.hdslb.com/bfs/openplatform/202403/u7J0bKWy1711434297039.png"
$ws1.Range("C26").Value = "杭州·第五人格ONLY"
$ws1.Range("D26").Value = "望江东路333号 瑞莱克斯大酒店"
$ws1.Range("E26").Value = "2024.04.20 10:00-04.20 17:00"
$ws1.Range("F26").Value = 686
$ws1.Range("G26").Value = 60
$ws1.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=81987"
$ws1.Range("I26").Value = "//i1.hdslb.com/bfs/openplatform/202402/Dxk0hWDb1708572766103.jpeg"
$ws1.Range("C27").Value = "杭州·黑执事only2.0"
$ws1.Range("D27").Value = "转塘街道创意路1号 艺创小镇凤凰创意大厦"
$ws1.Range("E27").Value = "2024.04.20 12:00-04.20 18:00"
$ws1.Range("F27").Value = 95
$ws1.Range("G27").Value = 98
$ws1.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=83170"
$ws1.Range("I27").Value = "//i1.hdslb.com/bfs/openplatform/202403/EZ57zs4Q1711004359139.jpeg"
$ws1.Range("B28").NumberFormat = "@"
$ws1.Range("B28").Value = "2024-04-30"
$ws1.Range("C28").Value = "杭州·造梦探险家——二次元同好会"
$ws1.Range("D28").Value = "临平街道北沙西路156-1号 杭州临平遇上设计师酒店"
$ws1.Range("E28").Value = "2024.04.30 10:00-05.01 16:00"
$ws1.Range("F28").Value = 76
$ws1.Range("G28").Value = 28
$ws1.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=82736"
$ws1.Range("I28").Value = "//i0.hdslb.com/bfs/openplatform/202403/RI8IceUW1710231050911.png"
$ws1.Range("B29").NumberFormat = "@"
$ws1.Range("B29").Value = "2024-05-01"
$ws1.Range("C29").Value = "杭州·与梦回望动漫游戏展"
$ws1.Range("D29").Value = "沈半路171号 T-Car杭州汽车文化主题公园"
$ws1.Range("E29").Value = "2024.05.01 10:00-05.02 17:00"
$ws1.Range("F29").Value = 391
$ws1.Range("G29").Value = 70
$ws1.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=82725"
$ws1.Range("I29").Value = "//i0.hdslb.com/bfs/openplatform/202403/lt13shal1710228931298.jpeg"
$ws1.Range("F30").Value = 391
$ws1.Range("F31").Value = 43
$ws1.Range("F33").Value = 969
$ws1.Range("F34").Value = 92
$ws1.Range("F35").Value = 13
$ws1.Range("F36").Value = 1061
$ws1.Range("F37").Value = 2013
$ws1.Range("F38").Value = 251
$ws1.Range("F39").Value = 9
$ws1.Range("F40").Value = 534
$ws1.Range("F43").Value = 636
$ws1.Range("F44").Value = 1292
$ws1.Range("F47").Value = 423

# ---- 演出 ----
$ws2.Range("G4").Value = 100

# ---- 全部类型 ----
$ws4.Range("F2").Value = 870
$ws4.Range("F3").Value = 1179
$ws4.Range("F5").Value = 65
$ws4.Range("F6").Value = 4309
$ws4.Range("F7").Value = 2578
$ws4.Range("F8").Value = 2477
$ws4.Range("F13").Value = 107
$ws4.Range("F16").Value = 269
$ws4.Range("F18").Value = 468
$ws4.Range("F21").Value = 523
$ws4.Range("F22").Value = 686
$ws4.Range("F23").Value = 95
$ws4.Range("G24").Value = 100
$ws4.Range("F27").Value = 76
$ws4.Range("F28").Value = 391
$ws4.Range("F30").Value = 969
$ws4.Range("F31").Value = 92
$ws4.Range("F33").Value = 1061
$ws4.Range("F34").Value = 2013
$ws4.Range("F35").Value = 251
$ws4.Range("F38").Value = 9
$ws4.Range("F40").Value = 534
$ws4.Range("F43").Value = 636
$ws4.Range("F44").Value = 1292
$ws4.Range("F47").Value = 423

